$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.127.71'
$ws.Range("E2").Value = '  +0.47%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.834.85'
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.41%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.58'
$ws.Range("E5").Value = '  +0.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6280'
$ws.Range("E6").Value = '  +0.30%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.003'
$ws.Range("E7").Value = '  +0.42%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07497'
$ws.Range("E8").Value = '  -1.19%  '
$ws.Range("E9").Value = '  +0.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.29'
$ws.Range("E10").Value = '  +3.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07690'
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.834.32'
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("E13").Value = '  +1.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6669'
$ws.Range("E14").Value = '  +0.45%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.82'
$ws.Range("E15").Value = '  +0.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009372'
$ws.Range("E16").Value = '  -8.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.980'
$ws.Range("E17").Value = '  -1.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.132.83'
$ws.Range("E18").Value = '  +0.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.081.95'
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("E20").Value = '  +2.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '223.29'
$ws.Range("E21").Value = '  -1.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.004'
$ws.Range("E22").Value = '  +0.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.096'
$ws.Range("E23").Value = '  -1.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.003'
$ws.Range("E24").Value = '  +0.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '160.08'
$ws.Range("E25").Value = '  +1.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1391'
$ws.Range("E26").Value = '  +1.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.499'
$ws.Range("E27").Value = '  +0.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.91'
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.502'
$ws.Range("E29").Value = '  +0.83%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05637'
$ws.Range("E30").Value = '  +7.90%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.157'
$ws.Range("E31").Value = '  +1.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.087'
$ws.Range("E32").Value = '  +2.00%  '
$ws.Range("E33").Value = '  +1.85%  '
$ws.Range("E34").Value = '  +0.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7421'
$ws.Range("E35").Value = '  +0.93%  '
$ws.Range("E36").Value = '  +0.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.674'
$ws.Range("E37").Value = '  -0.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.763'
$ws.Range("E38").Value = '  +0.29%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.221.78'
$ws.Range("E39").Value = '  -1.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01780'
$ws.Range("E40").Value = '  -0.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.545'
$ws.Range("E41").Value = '  +2.92%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8900'
$ws.Range("E42").Value = '  -0.77%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.003'
$ws.Range("E43").Value = '  +0.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.99'
$ws.Range("E44").Value = '  +0.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.980.91'
$ws.Range("E45").Value = '  -0.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.81'
$ws.Range("E46").Value = '  +2.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000123'
$ws.Range("E47").Value = '  -2.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5098'
$ws.Range("E48").Value = '  -0.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4076'
$ws.Range("E49").Value = '  +1.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07379'
$ws.Range("E50").Value = '  +5.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.006'
$ws.Range("E51").Value = '  +1.73%  '
